# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (interested-count) values in column F
# on sheets "展览" and "全部类型", and sheet "演出".

$wb = $excel.ActiveWorkbook

# Map: sheet name -> list of (row, newValue)
$changesBySheet = @{
    "展览" = @(
        @{Row = 3;  Value = 16517}
        @{Row = 4;  Value = 23}
        @{Row = 6;  Value = 15639}
        @{Row = 8;  Value = 511}
        @{Row = 10; Value = 107}
        @{Row = 11; Value = 653}
        @{Row = 15; Value = 1168}
        @{Row = 17; Value = 28}
        @{Row = 19; Value = 562}
        @{Row = 25; Value = 76}
        @{Row = 26; Value = 280}
        @{Row = 27; Value = 383}
        @{Row = 28; Value = 483}
        @{Row = 30; Value = 5831}
        @{Row = 31; Value = 5268}
    )
    "演出" = @(
        @{Row = 2;  Value = 84}
    )
    "全部类型" = @(
        @{Row = 3;  Value = 16517}
        @{Row = 4;  Value = 23}
        @{Row = 6;  Value = 15639}
        @{Row = 8;  Value = 511}
        @{Row = 10; Value = 107}
        @{Row = 11; Value = 653}
        @{Row = 15; Value = 1168}
        @{Row = 17; Value = 28}
        @{Row = 19; Value = 562}
        @{Row = 22; Value = 84}
        @{Row = 27; Value = 76}
        @{Row = 28; Value = 280}
        @{Row = 29; Value = 383}
        @{Row = 30; Value = 483}
        @{Row = 32; Value = 5831}
        @{Row = 34; Value = 5268}
    )
}

foreach ($sheetName in $changesBySheet.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($change in $changesBySheet[$sheetName]) {
        $ws.Range("F" + $change.Row).Value = $change.Value
    }
}
